$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a full row (A:G) as text values (matching the source
# workbook's convention of storing these figures as inline/shared text
# rather than numbers), without leaving a residual "@" text style behind.
function Set-RowText($r, $a, $b, $c, $d, $e, $f, $g) {
    $rng = $ws.Range("A$r`:G$r")
    $rng.NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value2 = $a
    $ws.Cells.Item($r, 2).Value2 = $b
    $ws.Cells.Item($r, 3).Value2 = $c
    $ws.Cells.Item($r, 4).Value2 = $d
    $ws.Cells.Item($r, 5).Value2 = $e
    $ws.Cells.Item($r, 6).Value2 = $f
    $ws.Cells.Item($r, 7).Value2 = $g
    $rng.Style = "Normal"
}

Set-RowText 2 "06/08/8000" "5001.00" "5001.00" "5050.00" "5050.00" "49.00" "100.98"
Set-RowText 3 "08/06/2000" "4000.00" "9001.00" "4000.00" "9050.00" "49.00" "100.54"
Set-RowText 4 "03/08/2023" "4000.00" "13001.00" "5000.00" "14050.00" "1049.00" "108.07"
Set-RowText 5 "03/08/2023" "4740.00" "17741.00" "4041.00" "18091.00" "350.00" "101.97"
Set-RowText 6 "03/08/2023" "4141.00" "21882.00" "4142.00" "22233.00" "351.00" "101.60"
Set-RowText 7 "03/08/2023" "4000.00" "25882.00" "4000.00" "26233.00" "351.00" "101.36"
Set-RowText 8 "05/08/2023" "8000.00" "33882.00" "8000.00" "34233.00" "351.00" "101.04"
